$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(4, 8).Value = 566.3333  # H4 (was 700)
$ws.Cells.Item(4, 9).Value = 566.3333  # I4 (was 700)
$ws.Cells.Item(4, 11).Value = 566.3333  # K4 (was 700)
$ws.Cells.Item(4, 13).Value = -452.3333  # M4 (was -586)
$ws.Cells.Item(12, 8).Value = 587.25  # H12 (was 591.125)
$ws.Cells.Item(12, 9).Value = 317.25  # I12 (was 279.8)
$ws.Cells.Item(12, 10).Value = 857.25  # J12 (was 1110)
$ws.Cells.Item(12, 11).Value = 317.25  # K12 (was 279.8)
$ws.Cells.Item(12, 12).Value = 857.25  # L12 (was 1110)
$ws.Cells.Item(12, 13).Value = -147.25  # M12 (was -109.8)
$ws.Cells.Item(12, 14).Value = -1197.25  # N12 (was -1450)
$ws.Cells.Item(19, 8).Value = 940.2857  # H19 (was 937)
$ws.Cells.Item(19, 10).Value = 919.93335  # J19 (was 916.6875)
$ws.Cells.Item(19, 12).Value = 919.93335  # L19 (was 916.6875)
$ws.Cells.Item(19, 14).Value = -1269.93335  # N19 (was -1266.6875)
$ws.Cells.Item(133, 8).Value = 106020.09  # H133 (was 115385)
$ws.Cells.Item(133, 10).Value = 106020.09  # J133 (was 115385)
$ws.Cells.Item(133, 12).Value = 106020.09  # L133 (was 115385)
$ws.Cells.Item(133, 14).Value = -116140.09  # N133 (was -125505)
$ws.Cells.Item(136, 8).Value = 0  # H136 (was 55390)
$ws.Cells.Item(136, 9).Value = 0  # I136 (was 50000)
$ws.Cells.Item(136, 10).Value = 0  # J136 (was 60780)
$ws.Cells.Item(136, 11).Value = 0  # K136 (was 50000)
$ws.Cells.Item(136, 12).Value = 0  # L136 (was 60780)
$ws.Cells.Item(136, 13).ClearContents()  # M136 (was -44900)
$ws.Cells.Item(136, 14).ClearContents()  # N136 (was -70980)
$ws.Cells.Item(138, 8).Value = 4132.325  # H138 (was 4080.9048)
$ws.Cells.Item(138, 9).Value = 801.7308  # I138 (was 839.2)
$ws.Cells.Item(138, 10).Value = 10317.714  # J138 (was 8848.117)
$ws.Cells.Item(138, 11).Value = 2405.1924  # K138 (was 2517.6)
$ws.Cells.Item(138, 12).Value = 30953.142  # L138 (was 26544.351)
$ws.Cells.Item(138, 13).Value = 2734.8076  # M138 (was 2622.4)
$ws.Cells.Item(138, 14).Value = -41233.142  # N138 (was -36824.351)
$ws.Cells.Item(140, 8).Value = 59666  # H140 (was 59999)
$ws.Cells.Item(140, 10).Value = 59666  # J140 (was 59999)
$ws.Cells.Item(140, 12).Value = 59666  # L140 (was 59999)
$ws.Cells.Item(140, 14).Value = -70026  # N140 (was -70359)

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 1787935.5  # H32 (was 1813847)
$ws.Cells.Item(32, 9).Value = 1840418.2  # I32 (was 1867886.6)
$ws.Cells.Item(32, 11).Value = 1840418.2  # K32 (was 1867886.6)
$ws.Cells.Item(32, 13).Value = -1840131.2  # M32 (was -1867599.6)
$ws.Cells.Item(45, 8).Value = 3435.1333  # H45 (was 3332.8125)
$ws.Cells.Item(45, 9).Value = 1914.625  # I45 (was 1901.6666)
$ws.Cells.Item(45, 11).Value = 1914.625  # K45 (was 1901.6666)
$ws.Cells.Item(45, 13).Value = -1537.625  # M45 (was -1524.6666)
$ws.Cells.Item(74, 8).Value = 35341.938  # H74 (was 36691.305)
$ws.Cells.Item(74, 10).Value = 4730.125  # J74 (was 4871.3335)
$ws.Cells.Item(74, 12).Value = 4730.125  # L74 (was 4871.3335)
$ws.Cells.Item(74, 14).Value = -6478.125  # N74 (was -6619.3335)
$ws.Cells.Item(77, 8).Value = 35341.938  # H77 (was 36691.305)
$ws.Cells.Item(77, 10).Value = 4730.125  # J77 (was 4871.3335)
$ws.Cells.Item(77, 12).Value = 23650.625  # L77 (was 24356.6675)
$ws.Cells.Item(77, 14).Value = -32386.625  # N77 (was -33092.6675)
$ws.Cells.Item(110, 8).Value = 1445.1111  # H110 (was 1004.13336)
$ws.Cells.Item(110, 9).Value = 1375.75  # I110 (was 933)
$ws.Cells.Item(110, 11).Value = 1375.75  # K110 (was 933)
$ws.Cells.Item(110, 13).Value = 669.25  # M110 (was 1112)
$ws.Cells.Item(122, 8).Value = 15999.474  # H122 (was 16955.723)
$ws.Cells.Item(122, 9).Value = 17999.4  # I122 (was 19371.715)
$ws.Cells.Item(122, 11).Value = 53998.2  # K122 (was 58115.145)
$ws.Cells.Item(122, 13).Value = -51548.2  # M122 (was -55665.145)

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(60, 8).Value = 92331  # H60 (was 106748.5)
$ws.Cells.Item(60, 10).Value = 92331  # J60 (was 106748.5)
$ws.Cells.Item(60, 12).Value = 92331  # L60 (was 106748.5)
$ws.Cells.Item(60, 14).Value = -93529  # N60 (was -107946.5)
$ws.Cells.Item(134, 8).Value = 3838.2876  # H134 (was 4018.7646)
$ws.Cells.Item(134, 9).Value = 2954.0178  # I134 (was 3122.6604)
$ws.Cells.Item(134, 10).Value = 6751.1763  # J134 (was 7185)
$ws.Cells.Item(134, 11).Value = 8862.053400000001  # K134 (was 9367.9812)
$ws.Cells.Item(134, 12).Value = 20253.5289  # L134 (was 21555)
$ws.Cells.Item(134, 13).Value = -6327.053400000001  # M134 (was -6832.9812)
$ws.Cells.Item(134, 14).Value = -25323.5289  # N134 (was -26625)

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 5083.6113  # H16 (was 5237.0527)
$ws.Cells.Item(16, 9).Value = 3227.818  # I16 (was 3182.3635)
$ws.Cells.Item(16, 10).Value = 7999.857  # J16 (was 8062.25)
$ws.Cells.Item(16, 11).Value = 3227.818  # K16 (was 3182.3635)
$ws.Cells.Item(16, 12).Value = 7999.857  # L16 (was 8062.25)
$ws.Cells.Item(16, 13).Value = -2940.818  # M16 (was -2895.3635)
$ws.Cells.Item(16, 14).Value = -8573.857  # N16 (was -8636.25)
$ws.Cells.Item(28, 8).Value = 0  # H28 (was 48000)
$ws.Cells.Item(28, 10).Value = 0  # J28 (was 48000)
$ws.Cells.Item(28, 12).Value = 0  # L28 (was 48000)
$ws.Cells.Item(28, 14).ClearContents()  # N28 (was -48490)
$ws.Cells.Item(31, 8).Value = 7961.5  # H31 (was 8273.069)
$ws.Cells.Item(31, 9).Value = 5632.125  # I31 (was 6041.6)
$ws.Cells.Item(31, 11).Value = 5632.125  # K31 (was 6041.6)
$ws.Cells.Item(31, 13).Value = -5337.125  # M31 (was -5746.6)
$ws.Cells.Item(34, 8).Value = 7961.5  # H34 (was 8273.069)
$ws.Cells.Item(34, 9).Value = 5632.125  # I34 (was 6041.6)
$ws.Cells.Item(34, 11).Value = 5632.125  # K34 (was 6041.6)
$ws.Cells.Item(34, 13).Value = -5430.125  # M34 (was -5839.6)
$ws.Cells.Item(58, 8).Value = 10642405  # H58 (was 10873661)
$ws.Cells.Item(58, 9).Value = 17243328  # I58 (was 17243320)
$ws.Cells.Item(58, 10).Value = 7585.8887  # J58 (was 7770.7646)
$ws.Cells.Item(58, 11).Value = 17243328  # K58 (was 17243320)
$ws.Cells.Item(58, 12).Value = 7585.8887  # L58 (was 7770.7646)
$ws.Cells.Item(58, 13).Value = -17243125  # M58 (was -17243117)
$ws.Cells.Item(58, 14).Value = -7991.8887  # N58 (was -8176.7646)
$ws.Cells.Item(113, 8).Value = 5083.6113  # H113 (was 5237.0527)
$ws.Cells.Item(113, 9).Value = 3227.818  # I113 (was 3182.3635)
$ws.Cells.Item(113, 10).Value = 7999.857  # J113 (was 8062.25)
$ws.Cells.Item(113, 11).Value = 3227.818  # K113 (was 3182.3635)
$ws.Cells.Item(113, 12).Value = 7999.857  # L113 (was 8062.25)
$ws.Cells.Item(113, 13).Value = -1057.818  # M113 (was -1012.3635)
$ws.Cells.Item(113, 14).Value = -12339.857  # N113 (was -12402.25)
$ws.Cells.Item(132, 8).Value = 4082.577  # H132 (was 4112.923)
$ws.Cells.Item(132, 9).Value = 1861.8823  # I132 (was 1908.2941)
$ws.Cells.Item(132, 11).Value = 5585.6469  # K132 (was 5724.8823)
$ws.Cells.Item(132, 13).Value = -3055.6469  # M132 (was -3194.8823)
$ws.Cells.Item(134, 8).Value = 6877.3096  # H134 (was 7505.108)
$ws.Cells.Item(134, 9).Value = 6768.6333  # I134 (was 7418.4814)
$ws.Cells.Item(134, 10).Value = 7149  # J134 (was 7739)
$ws.Cells.Item(134, 11).Value = 20305.8999  # K134 (was 22255.4442)
$ws.Cells.Item(134, 12).Value = 21447  # L134 (was 23217)
$ws.Cells.Item(134, 13).Value = -17770.8999  # M134 (was -19720.4442)
$ws.Cells.Item(134, 14).Value = -26517  # N134 (was -28287)
$ws.Cells.Item(136, 8).Value = 10642405  # H136 (was 10873661)
$ws.Cells.Item(136, 9).Value = 17243328  # I136 (was 17243320)
$ws.Cells.Item(136, 10).Value = 7585.8887  # J136 (was 7770.7646)
$ws.Cells.Item(136, 11).Value = 51729984  # K136 (was 51729960)
$ws.Cells.Item(136, 12).Value = 22757.6661  # L136 (was 23312.2938)
$ws.Cells.Item(136, 13).Value = -51727434  # M136 (was -51727410)
$ws.Cells.Item(136, 14).Value = -27857.6661  # N136 (was -28412.2938)

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(14, 8).Value = 20833574  # H14 (was 33333550)
$ws.Cells.Item(14, 9).Value = 20833574  # I14 (was 33333550)
$ws.Cells.Item(14, 11).Value = 62500722  # K14 (was 100000650)
$ws.Cells.Item(14, 13).Value = -62500549  # M14 (was -100000477)
$ws.Cells.Item(97, 8).Value = 151.33333  # H97 (was 374.5)
$ws.Cells.Item(97, 9).Value = 200  # I97 (was 500)
$ws.Cells.Item(97, 10).Value = 127  # J97 (was 249)
$ws.Cells.Item(97, 11).Value = 600  # K97 (was 1500)
$ws.Cells.Item(97, 12).Value = 381  # L97 (was 747)
$ws.Cells.Item(97, 13).Value = -104  # M97 (was -1004)
$ws.Cells.Item(97, 14).Value = -1373  # N97 (was -1739)

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 102.3  # H2 (was 103.8)
$ws.Cells.Item(2, 9).Value = 58.11111  # I2 (was 59.77778)
$ws.Cells.Item(2, 11).Value = 58.11111  # K2 (was 59.77778)
$ws.Cells.Item(2, 13).Value = 54.88889  # M2 (was 53.22222)
$ws.Cells.Item(126, 8).Value = 3139.7  # H126 (was 2924.5833)
$ws.Cells.Item(126, 9).Value = 2682.8333  # I126 (was 2683)
$ws.Cells.Item(126, 10).Value = 3825  # J126 (was 3166.1667)
$ws.Cells.Item(126, 11).Value = 8048.499899999999  # K126 (was 8049)
$ws.Cells.Item(126, 12).Value = 11475  # L126 (was 9498.500100000001)
$ws.Cells.Item(126, 13).Value = -5578.499899999999  # M126 (was -5579)
$ws.Cells.Item(126, 14).Value = -16415  # N126 (was -14438.5001)
$ws.Cells.Item(132, 8).Value = 2622.68  # H132 (was 2645.6458)
$ws.Cells.Item(132, 9).Value = 2576.756  # I132 (was 2623.7)
$ws.Cells.Item(132, 10).Value = 2831.889  # J132 (was 2755.375)
$ws.Cells.Item(132, 11).Value = 7730.268  # K132 (was 7871.099999999999)
$ws.Cells.Item(132, 12).Value = 8495.667000000001  # L132 (was 8266.125)
$ws.Cells.Item(132, 13).Value = -5200.268  # M132 (was -5341.099999999999)
$ws.Cells.Item(132, 14).Value = -13555.667  # N132 (was -13326.125)

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(122, 8).Value = 4528.3477  # H122 (was 4637.3022)
$ws.Cells.Item(122, 9).Value = 3996.7666  # I122 (was 4107.7407)
$ws.Cells.Item(122, 10).Value = 5525.0625  # J122 (was 5530.9375)
$ws.Cells.Item(122, 11).Value = 11990.2998  # K122 (was 12323.2221)
$ws.Cells.Item(122, 12).Value = 16575.1875  # L122 (was 16592.8125)
$ws.Cells.Item(122, 13).Value = -9540.299800000001  # M122 (was -9873.222100000001)
$ws.Cells.Item(122, 14).Value = -21475.1875  # N122 (was -21492.8125)
$ws.Cells.Item(132, 8).Value = 10875560  # H132 (was 11911139)
$ws.Cells.Item(132, 9).Value = 21744536  # I132 (was 26321916)
$ws.Cells.Item(132, 11).Value = 65233608  # K132 (was 78965748)
$ws.Cells.Item(132, 13).Value = -65231078  # M132 (was -78963218)
$ws.Cells.Item(136, 8).Value = 8128.4707  # H136 (was 8206.714)
$ws.Cells.Item(136, 9).Value = 7836.8  # I136 (was 7947)
$ws.Cells.Item(136, 11).Value = 23510.4  # K136 (was 23841)
$ws.Cells.Item(136, 13).Value = -20960.4  # M136 (was -21291)

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(13, 8).Value = 300004320  # H13 (was 300005340)
$ws.Cells.Item(13, 10).Value = 450002500  # J13 (was 900000000)
$ws.Cells.Item(13, 12).Value = 450002500  # L13 (was 900000000)
$ws.Cells.Item(13, 14).Value = -450002780  # N13 (was -900000280)
$ws.Cells.Item(122, 8).Value = 238052.94  # H122 (was 269713.34)
$ws.Cells.Item(122, 9).Value = 446356.12  # I122 (was 573629.3)
$ws.Cells.Item(122, 10).Value = 3711.875  # J122 (was 3786.875)
$ws.Cells.Item(122, 11).Value = 1339068.36  # K122 (was 1720887.9)
$ws.Cells.Item(122, 12).Value = 11135.625  # L122 (was 11360.625)
$ws.Cells.Item(122, 13).Value = -1336618.36  # M122 (was -1718437.9)
$ws.Cells.Item(122, 14).Value = -16035.625  # N122 (was -16260.625)
$ws.Cells.Item(126, 8).Value = 4900  # H126 (was 6275)
$ws.Cells.Item(126, 9).Value = 4300  # I126 (was 4550)
$ws.Cells.Item(126, 10).Value = 5700  # J126 (was 8000)
$ws.Cells.Item(126, 11).Value = 12900  # K126 (was 13650)
$ws.Cells.Item(126, 12).Value = 17100  # L126 (was 24000)
$ws.Cells.Item(126, 13).Value = -10430  # M126 (was -11180)
$ws.Cells.Item(126, 14).Value = -22040  # N126 (was -28940)
$ws.Cells.Item(132, 8).Value = 12509227  # H132 (was 12509172)
$ws.Cells.Item(132, 9).Value = 14708738  # I132 (was 14288511)
$ws.Cells.Item(132, 10).Value = 45331.332  # J132 (was 53798.6)
$ws.Cells.Item(132, 11).Value = 44126214  # K132 (was 42865533)
$ws.Cells.Item(132, 12).Value = 135993.996  # L132 (was 161395.8)
$ws.Cells.Item(132, 13).Value = -44123684  # M132 (was -42863003)
$ws.Cells.Item(132, 14).Value = -141053.996  # N132 (was -166455.8)
$ws.Cells.Item(136, 8).Value = 50054690  # H136 (was 47671180)
$ws.Cells.Item(136, 9).Value = 100001740  # I136 (was 90910770)
$ws.Cells.Item(136, 11).Value = 300005220  # K136 (was 272732310)
$ws.Cells.Item(136, 13).Value = -300002670  # M136 (was -272729760)
